# Applies the "deploying to gh-pages" content update to the FHIR
# StructureDefinition workbook:
#  - Metadata sheet: bump Version, Date; fill in Publisher; replace the
#    duplicated "Contact" rows with a single "Jurisdiction" row.
#  - Elements sheet: give the root Extension row a specific Short/Definition
#    (matching the resource's own Title/Description) instead of the generic
#    "Extension" / "An Extension" placeholder text.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value2 = "6.0.0"

# Date: updated publication timestamp
$meta.Range("B8").Value2 = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated
$meta.Range("B9").Value2 = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail"; repurpose it as
# the new "Jurisdiction" row
$meta.Range("A10").Value2 = "Jurisdiction"
$meta.Range("B10").Value2 = "United States of America"

# Row 11 duplicated the old "Contact" row and is no longer needed; remove it
# so every row below shifts up by one (A1:B21 -> A1:B20)
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

# Root Extension row (row 2): Short/Definition columns (K/L) get the
# resource-specific text instead of the generic placeholders
$elements.Range("K2").Value2 = "Employee Wage Basis"
$elements.Range("L2").Value2 = "Code indicating the basis on which the wages of the employee are calculated (e.g., hourly, daily, weekly, bimonthly, monthly, annually)"
